$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("C1").Value = "rules"
$ws.Range("E1").Value = "adaptive_filter"

# Values for E (adaptive_filter name), F (RMSE), G (NDEI), H (MAE) for rows 2-13
$data = @{
    2  = @{ E = "RLS"; F = [double]"4.714301534104477e-14"; G = [double]"5.758793880241753e-15"; H = [double]"3.543533695637979e-14" }
    3  = @{ E = "RLS"; F = [double]"1.12536913539333e-14";  G = [double]"1.374703939286129e-15"; H = [double]"8.907168405625221e-15" }
    4  = @{ E = "RLS"; F = [double]"1.087463607170333e-14"; G = [double]"1.328400128980666e-15"; H = [double]"9.715101986773611e-15" }
    5  = @{ E = "RLS"; F = [double]"1.005715688141166e-14"; G = [double]"1.228540284967295e-15"; H = [double]"7.926785963729976e-15" }
    6  = @{ E = "RLS"; F = [double]"5.865903277367848e-15"; G = [double]"7.165542477802689e-16"; H = [double]"4.557647662051245e-15" }
    7  = @{ E = "RLS"; F = [double]"5.606910118278322e-15"; G = [double]"6.849167250465314e-16"; H = [double]"4.675799678000025e-15" }
    8  = @{ E = "RLS"; F = [double]"1.318412360474831e-14"; G = [double]"1.610517481372731e-15"; H = [double]"1.129639264074722e-14" }
    9  = @{ E = "RLS"; F = [double]"5.458280528191767e-15"; G = [double]"6.667607550131915e-16"; H = [double]"4.378582565967015e-15" }
    10 = @{ E = "RLS"; F = [double]"5.15702267032874e-15";  G = [double]"6.299603531787785e-16"; H = [double]"3.644836168992249e-15" }
    11 = @{ E = "RLS"; F = [double]"1.687181030391637e-14"; G = [double]"2.06098988840454e-15";  H = [double]"1.477405524308306e-14" }
    12 = @{ E = "RLS"; F = [double]"1.565610876326842e-14"; G = [double]"1.912484865086938e-15"; H = [double]"1.47312318593551e-14" }
    13 = @{ E = "RLS"; F = [double]"3.952490733747522e-05"; G = [double]"4.828197620486155e-06"; H = [double]"2.950236822960124e-05" }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
}
